# Update "想去人数" (column F) counts across all sheets to the latest
# scraped values, as published to gh-pages at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览" = @{
        2 = 2902; 3 = 21329; 5 = 3174; 6 = 821; 8 = 529; 9 = 785; 10 = 290;
        11 = 266; 13 = 128; 14 = 537; 15 = 186; 16 = 293; 19 = 98; 20 = 33;
        21 = 27; 22 = 48; 23 = 135
    }
    "演出" = @{
        2 = 35; 4 = 349; 11 = 103; 13 = 161
    }
    "本地生活" = @{
        2 = 6147; 3 = 714; 4 = 708; 5 = 1675; 6 = 64
    }
    "全部类型" = @{
        2 = 6147; 3 = 714; 4 = 708; 5 = 1675; 6 = 2902; 7 = 21329; 8 = 35;
        11 = 349; 12 = 3174; 13 = 821; 15 = 64; 17 = 529; 18 = 785; 19 = 290;
        20 = 266; 25 = 128; 28 = 537; 29 = 103; 30 = 186; 32 = 293; 33 = 161;
        34 = 161; 38 = 98; 39 = 33; 42 = 27; 43 = 48; 49 = 135
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowsMap = $updates[$sheetName]
    foreach ($row in $rowsMap.Keys) {
        $ws.Cells.Item($row, 6).Value = $rowsMap[$row]
    }
}
